$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6804673671722412
$ws.Range("B1").Value = 0.5094919204711914
$ws.Range("C1").Value = 0.5514360070228577
$ws.Range("D1").Value = 3.89283561706543
$ws.Range("E1").Value = 1.591897964477539
